$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.6967741935483871
$ws.Range("C2").Value = 0.7448275862068966
$ws.Range("D2").Value = 0.7200000000000001
$ws.Range("E2").Value = 145

# Row 3
$ws.Range("B3").Value = 0.8472222222222222
$ws.Range("C3").Value = 0.8243243243243243
$ws.Range("D3").Value = 0.8356164383561645
$ws.Range("E3").Value = 148

# Row 4
$ws.Range("B4").Value = 0.7435897435897436
$ws.Range("C4").Value = 0.7682119205298014
$ws.Range("D4").Value = 0.7557003257328991
$ws.Range("E4").Value = 151

# Row 5
$ws.Range("B5").Value = 0.6896551724137931
$ws.Range("C5").Value = 0.6410256410256411
$ws.Range("D5").Value = 0.6644518272425251
$ws.Range("E5").Value = 156

# Row 6 (accuracy row)
$ws.Range("B6").Value = 0.7433333333333333
$ws.Range("C6").Value = 0.7433333333333333
$ws.Range("D6").Value = 0.7433333333333333
$ws.Range("E6").Value = 0.7433333333333333

# Row 7 (macro avg)
$ws.Range("B7").Value = 0.7443103329435365
$ws.Range("C7").Value = 0.7445973680216658
$ws.Range("D7").Value = 0.7439421478328971

# Row 8 (weighted avg)
$ws.Range("B8").Value = 0.7438156752200134
$ws.Range("C8").Value = 0.7433333333333333
$ws.Range("D8").Value = 0.7430607785203568
